$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 1.9
$ws.Range("H3").Value = 3.2
$ws.Range("I3").Value = 4.75
$ws.Range("J3").Value = 2.63
$ws.Range("L3").Value = 5.5
$ws.Range("Q3").Value = 1.93
$ws.Range("R3").Value = 1.93
$ws.Range("S3").Value = 2.5
$ws.Range("T3").Value = 1.5
$ws.Range("U3").Value = 4.2
$ws.Range("V3").Value = 1.23
$ws.Range("AA3").Value = 2.2
$ws.Range("AB3").Value = 1.62
$ws.Range("AD3").Value = 7.5
$ws.Range("AF3").Value = 15
$ws.Range("AJ3").Value = 6.5
$ws.Range("AK3").Value = 21
$ws.Range("AL3").Value = 81
$ws.Range("AP3").Value = 17
